$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 25 is no longer the last row of the sheet: it becomes the final
# row of its data block, so it picks up the "block end" bottom-border
# formatting (same look as row 17 / row 21). Copy formats only (values
# are untouched) from row 21, which already has that formatting.
$ws.Range("A21:E21").Copy()
$ws.Range("A25").PasteSpecial(-4122)

# --- New rows 26 & 27: values are written in the same order the source
# file's string table was built in, so new shared-string entries land at
# the same indices as the target workbook.
$ws.Range("C26").Value = " We owe a lot to you."
$ws.Range("C27").Value = " I can relax like this because it\'s\npeaceful.[K] Thanks to you!"
$ws.Range("A26").Value = "SCRIPT/P02P01A/us0103.ssb "
$ws.Range("D26").Value = " Мы у вас в долгу."
$ws.Range("D27").Value = " В округе царит мир и я могу\nрасслабиться.[K] Всё благодаря вам!"
$ws.Range("E26").Value = " Íú ô âàò â äïìãô."
$ws.Range("E27").Value = " Â ïëñôãå øàñéó íéñ é ÿ íïãô\nñàòòìàáéóûòÿ.[K] Âòæ áìàãïäàñÿ âàí!"
$ws.Range("A27").Value = "SCRIPT/P02P01A/us3103.ssb"
$ws.Range("B26").Value = 18
$ws.Range("B27").Value = 21

# Row heights for the two new rows (matches the wrapped-text rows above).
$ws.Rows.Item(26).RowHeight = 43.2
$ws.Rows.Item(27).RowHeight = 43.2

# --- Move the selection down to the newly added last cell, mirroring
# where the author's cursor ended up after the edit.
$ws.Range("E27").Select()
